$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of "Kevin" test results, added under friends' advice
$rowValues = @{
    "B6" = "Yes"
    "C6" = "Yes"
    "D6" = "Yes"
    "E6" = "No, 9.39"
    "F6" = "Yes"
    "G6" = "Yes"
    "H6" = "No, 9.89"
    "I6" = "Yes"
    "J6" = "No, 9.89"
    "K6" = "Yes"
    "L6" = "No, 9.89"
    "M6" = "No, 10.39"
    "N6" = "No, 10.39"
    "O6" = "Yes"
    "P6" = "No, 10.39"
    "Q6" = "No, 10.39"
}

foreach ($addr in $rowValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $rowValues[$addr]
    $cell.HorizontalAlignment = -4108  # xlCenter
}

$ws.Range("M10").Select()
